$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testGoogle")

# --- Row 4: TC_003 / TRUE / TRUE -----------------------------------------
# Column A keeps the same bordered style used by the existing TC_001/TC_002
# rows (row 3) -- copy that formatting across before writing the value.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "TC_003"

# Columns B/C on the new rows were typed with a leading apostrophe (quote
# prefix) so Excel keeps "TRUE"/"10000100"/"10002" as literal text instead
# of coercing to boolean/number -- build that style once...
$ws.Range("B4").Value = "'TRUE"
$ws.Range("B4").NumberFormat = "@"
# ...then stamp the resulting (quote-prefixed, text-formatted) style onto
# the other three cells before filling them in, so every cell shares one
# single style entry instead of minting a new one each time.
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("C4").Value = "'TRUE"

# --- Row 5: TC_003 / 10000100 / 10002 ------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "TC_003"

$ws.Range("B5").Value = "'10000100"
$ws.Range("C5").Value = "'10002"

# Column B widens to fit the new longest entry ("10000100").
$ws.Columns.Item(2).AutoFit()

# testGoogle (was tab 2 / Sheet1) becomes the selected/active sheet.
$ws.Activate()
